$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Update existing field-layout values (row numbers shifted because new field inserted) ---
$ws.Range("B5").Value = 4
$ws.Range("D5").Value = 40

$ws.Range("C6").Value = 41
$ws.Range("D6").Value = 44

$ws.Range("C7").Value = 45
$ws.Range("D7").Value = 45

$ws.Range("C8").Value = 46
$ws.Range("D8").Value = 53

$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 54
$ws.Range("D9").Value = 55

$ws.Range("C10").Value = 56
$ws.Range("D10").Value = 58

# --- Add new row 11: MODREALACT field, duplicating row 10's formatting ---
$ws.Rows.Item(10).Copy()
$ws.Rows.Item(11).Insert(-4121)

$ws.Range("A11").Value = "MODREALACT"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 59
$ws.Range("D11").Value = 59
$ws.Range("E11").Value = "c"
$ws.Range("F11").Value = "Modalité de réalisation de l'acte"

# --- Column F widened to fit the new, longer label ---
$ws.Columns.Item(6).ColumnWidth = 26.59

# --- Selection moved to G4 (matches the saved view state) ---
$ws.Range("G4").Select() | Out-Null
